# Wrong calculation of total elec capacity for HTSE modules: include ANR
# Thermal efficiency (column K, "ANR Th Eff") in the HTSE electric-capacity
# formula; add the new column for all tech rows and fix up the PEM/Alkaline
# D-column shared-formula anchor that shifted down a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- New column header ---------------------------------------------------
$ws.Range("K1").Value = "ANR Th Eff"

# --- HTSE rows (2-6): pull ANR thermal efficiency into the elec-capacity calc
$ws.Range("C2").Formula = "=D2*(E2+(F2*K2))"
$ws.Range("C3").Formula = "=D3*(E3+(F3*K3))"
$ws.Range("C4").Formula = "=D4*(E4+(F4*K4))"
$ws.Range("C5").Formula = "=D5*(E5+(F5*K5))"
$ws.Range("C6").Formula = "=D6*(E6+(F6*K6))"

# --- New "ANR Th Eff" values for every tech/reactor-type row -------------
# HTGR reactor-type ordering across the three tech blocks: iPWR, HTGR,
# PBR-HTGR, iMSR, Micro -> 0.31, 0.47, 0.4, 0.47, 0.33
$kValues = @{
    2  = 0.31; 3  = 0.47; 4  = 0.4; 5  = 0.47; 6  = 0.33;
    7  = 0.31; 8  = 0.47; 9  = 0.4; 10 = 0.47; 11 = 0.33;
    12 = 0.31; 13 = 0.47; 14 = 0.4; 15 = 0.47; 16 = 0.33
}

foreach ($row in 2..16) {
    $cell = $ws.Range("K$row")
    $cell.Value = $kValues[$row]
    $cell.NumberFormat = "#,##0.00"
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Color = 0
    $cell.HorizontalAlignment = -4152   # xlRight
}

# --- Corrected thermal-eff HTGR/PBR shared formula in the Alkaline block -
# D14 used to anchor the D14:D16 shared formula; move the anchor down to
# D15 and give D14 its own standalone formula.
$ws.Range("D14").Formula = "=C14/E14"
$ws.Range("D15").Formula = "=C15/E15"
$ws.Range("D16").Formula = "=C16/E16"

# --- Leave the selection where the author finished editing ---------------
$ws.Activate()
$ws.Range("F32").Select()
